$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("D40").Value = "19/6/2025"
$ws.Range("E40").Value = 373
$ws.Range("F40").Value = 553
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 1012

$ws.Range("D41").Value = "20/6/2025"

$ws.Range("J40").Value = "No he podido adelantar, le dare en fin de semana, disculpa (rafael)"

$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J38").Select()
